# PetRegistry: add pet log rows + make filters work on the full data set.
# Target layout (A2:I9):
#   A: record Id, B: Species, C: Gender, D: Birth date, E/F: constants,
#   G: pet Name, H: owner comment ("-" when present), I: Owner.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = 2
$ws.Cells.Item(2,2).Value = "Кошка"
$ws.Cells.Item(2,3).Value = "Женский"
$ws.Cells.Item(2,4).Value = "19.09.2019"
$ws.Cells.Item(2,5).Value = 123
$ws.Cells.Item(2,6).Value = 123
$ws.Cells.Item(2,7).Value = "Вася"
$ws.Cells.Item(2,8).Value = $null
$ws.Cells.Item(2,9).Value = "Sapov EV"
# Row 3
$ws.Cells.Item(3,1).Value = 3
$ws.Cells.Item(3,2).Value = "Собака"
$ws.Cells.Item(3,3).Value = "Женский"
$ws.Cells.Item(3,4).Value = "19.09.2019"
$ws.Cells.Item(3,5).Value = 123
$ws.Cells.Item(3,6).Value = 123
$ws.Cells.Item(3,7).Value = "Тося"
$ws.Cells.Item(3,8).Value = "-"
$ws.Cells.Item(3,9).Value = "Ромашка"
# Row 4
$ws.Cells.Item(4,1).Value = 4
$ws.Cells.Item(4,2).Value = "Собака"
$ws.Cells.Item(4,3).Value = "Мужской"
$ws.Cells.Item(4,4).Value = "19.09.2019"
$ws.Cells.Item(4,5).Value = 123
$ws.Cells.Item(4,6).Value = 123
$ws.Cells.Item(4,7).Value = "Бося"
$ws.Cells.Item(4,8).Value = "-"
$ws.Cells.Item(4,9).Value = "Sapov EV"
# Row 5
$ws.Cells.Item(5,1).Value = 1002
$ws.Cells.Item(5,2).Value = "Кошка"
$ws.Cells.Item(5,3).Value = "Женский"
$ws.Cells.Item(5,4).Value = "19.09.2019"
$ws.Cells.Item(5,5).Value = 123
$ws.Cells.Item(5,6).Value = 123
$ws.Cells.Item(5,7).Value = "Жося"
$ws.Cells.Item(5,8).Value = "-"
$ws.Cells.Item(5,9).Value = "Sapov EV"
# Row 6
$ws.Cells.Item(6,1).Value = 1003
$ws.Cells.Item(6,2).Value = "Кошка"
$ws.Cells.Item(6,3).Value = "Мужской"
$ws.Cells.Item(6,4).Value = "19.09.2019"
$ws.Cells.Item(6,5).Value = 123
$ws.Cells.Item(6,6).Value = 123
$ws.Cells.Item(6,7).Value = "Гусь"
$ws.Cells.Item(6,8).Value = "-"
$ws.Cells.Item(6,9).Value = "Цветок"
# Row 7
$ws.Cells.Item(7,1).Value = 1004
$ws.Cells.Item(7,2).Value = "Собака"
$ws.Cells.Item(7,3).Value = "Мужской"
$ws.Cells.Item(7,4).Value = "19.09.2019"
$ws.Cells.Item(7,5).Value = 123
$ws.Cells.Item(7,6).Value = 123
$ws.Cells.Item(7,7).Value = "Дося"
$ws.Cells.Item(7,8).Value = "-"
$ws.Cells.Item(7,9).Value = "Lap YD"
# Row 8
$ws.Cells.Item(8,1).Value = 2002
$ws.Cells.Item(8,2).Value = "Собака"
$ws.Cells.Item(8,3).Value = "Мужской"
$ws.Cells.Item(8,4).Value = "04.06.2020"
$ws.Cells.Item(8,5).Value = 123
$ws.Cells.Item(8,6).Value = 123
$ws.Cells.Item(8,7).Value = "Жук"
$ws.Cells.Item(8,8).Value = $null
$ws.Cells.Item(8,9).Value = "Sapov EV"
# Row 9
$ws.Cells.Item(9,1).Value = 2007
$ws.Cells.Item(9,2).Value = "Собака"
$ws.Cells.Item(9,3).Value = "Женский"
$ws.Cells.Item(9,4).Value = "01.12.2021"
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = "Зая"
$ws.Cells.Item(9,8).Value = $null
$ws.Cells.Item(9,9).Value = "Ромашка"
